$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the INDEX field to INDEX_IN_MARKUP (row 11, column B)
$ws.Range("B11").Value = "INDEX_IN_MARKUP"

# New column E: SQL header + CREATE TABLE statements for the two tables
$ws.Range("E1").Value = "SQL"

$ws.Range("E5").Value = 'CREATE TABLE `spinnis0_WPLXP`.`_LXP_auto_feed` ( `ID` BIGINT UNSIGNED NOT NULL AUTO_INCREMENT , `FEED` TEXT CHARACTER SET utf8mb4 COLLATE utf8mb4_unicode_520_ci NOT NULL , PRIMARY KEY (`ID`)) ENGINE = MyISAM CHARSET=utf8mb4 COLLATE utf8mb4_unicode_ci;'

$ws.Range("E15").Value = 'CREATE TABLE `spinnis0_WPLXP`.`LXP_auto_feed_posted` ( `AUTO_FEED_ID` BIGINT UNSIGNED NOT NULL , `POST_ID` BIGINT UNSIGNED NOT NULL , `INDEX_IN_MARKUP` INT UNSIGNED NOT NULL , `TITLE` VARCHAR(512) NOT NULL , `ORIGINAL_POST_HYPERLINK` VARCHAR(512) NOT NULL , `DATE_CREATED` DATETIME NOT NULL DEFAULT CURRENT_TIMESTAMP , `IS_EDITED` TINYINT(1) NOT NULL DEFAULT ''0'' , `NOTES` VARCHAR(512) NULL ) ENGINE = MyISAM;'

# Update the print/page setup (orientation) for the worksheet
$ws.PageSetup.Orientation = 1

# Move the active selection to C26, as recorded in the saved view state
$ws.Range("C26").Select() | Out-Null
